# EmissionsTable.xlsx edit
#
# 1) "Emissions By Atomic Number": every placeholder cell that previously
#    held the em-dash ("-") lookup string now shows "N/A" instead. This
#    introduces a new shared string ("N/A") and repoints ~220 cells at it.
# 2) View/selection state: sheet 1 becomes the active tab (scrolled so the
#    data rows are visible) with F23 selected; sheet 2 keeps A1 scrolled
#    further down with D299 selected; sheet 3 is no longer the active tab,
#    with A17 selected.

$wb = $excel.ActiveWorkbook

$wsAtomic   = $wb.Worksheets.Item("Emissions By Atomic Number")
$wsAll      = $wb.Worksheets.Item("Emissions By Energy (All)")
$wsFiltered = $wb.Worksheets.Item("Emissions By Energy (Filtered)")

# --- 1) Replace every "-" placeholder cell with "N/A" on the atomic-number sheet ---
$naCells = "D3","E3","F3","G3","H3","I3","J3","K3","D4","E4","F4","G4","H4","I4","J4","K4","D5","E5","F5","G5","H5","I5","J5","K5","D6","E6","F6","G6","H6","I6","J6","K6","D7","E7","F7","G7","H7","I7","J7","K7","D8","E8","F8","G8","H8","I8","J8","K8","D9","E9","F9","G9","H9","I9","J9","K9","E10","F10","G10","H10","I10","J10","K10","F11","G11","H11","I11","J11","K11","F12","G12","H12","I12","J12","K12","F13","G13","H13","I13","J13","K13","F14","G14","H14","I14","J14","K14","F15","G15","H15","I15","J15","K15","F16","G16","H16","I16","J16","K16","F17","G17","H17","I17","J17","K17","F18","G18","H18","I18","J18","K18","F19","G19","H19","I19","J19","K19","I20","J20","K20","I21","J21","K21","I22","J22","K22","I23","J23","K23","I24","J24","K24","I25","J25","K25","I26","J26","K26","I27","J27","K27","I28","J28","K28","I29","J29","K29","I30","J30","K30","I31","J31","K31","I32","J32","K32","I33","J33","K33","I34","J34","K34","I35","J35","K35","I36","J36","K36","I37","J37","K37","I38","J38","K38","I39","J39","K39","K40","K41","K42","K43","K44","K45","K46","K47","K48","K49","K50","K51","K52","K53","G54","H54","I54","J54","K54","K55","K56","K61","K84","I85","K85","I86","K86","K87","K88","I89","K89","C93","D93","E93","K93","C94","D94","E94","K94","C95","D95","E95","K95"

foreach ($addr in $naCells) {
    $wsAtomic.Range($addr).Value = "N/A"
}

# --- 2) View / selection state ---
$wsAll.Activate()
$wsAll.Range("D299").Select()

$wsFiltered.Activate()
$wsFiltered.Range("A17").Select()

$wsAtomic.Activate()
$wsAtomic.Range("F23").Select()

Write-Output "done"
